$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the long-form metric header labels (E1:BN1) to their abbreviated forms.
# (The underlying shared-string table is re-packed by the engine on save; the
# old long labels are dropped automatically once no cell references them.)
$ws.Range("E1").Value = "BS"
$ws.Range("F1").Value = "Run_Env"
$ws.Range("G1").Value = "CPU_Util_%_ran"
$ws.Range("H1").Value = "Max_CPU_Util_%"
$ws.Range("I1").Value = "Min_CPU_Util_%"
$ws.Range("J1").Value = "GPU_Util_%_ran"
$ws.Range("K1").Value = "MGUP"
$ws.Range("L1").Value = "Min_GPU_Util_%"
$ws.Range("M1").Value = "Sys_Mem_Util_%_ran"
$ws.Range("N1").Value = "Max_Sys_Mem_Util_%"
$ws.Range("O1").Value = "Min_Sys_Mem_Util_%"
$ws.Range("P1").Value = "MaPMiU(non swap)"
$ws.Range("Q1").Value = "PMiU(non swap)_%"
$ws.Range("R1").Value = "CPU_Thds"
$ws.Range("S1").Value = "GPU_Temp_ran"
$ws.Range("T1").Value = "Max_GPU_Temp"
$ws.Range("U1").Value = "Min_GPU_Temp"
$ws.Range("V1").Value = "GTSAMPR"
$ws.Range("W1").Value = "MaGTSAMP"
$ws.Range("X1").Value = "MiGTSAMP"
$ws.Range("Y1").Value = "MGMAP"
$ws.Range("Z1").Value = "GPUPR"
$ws.Range("AA1").Value = "MaGPUPR"
$ws.Range("AB1").Value = "MiGPUPR"
$ws.Range("AC1").Value = "SCT_user"
$ws.Range("AD1").Value = "SCT_nice"
$ws.Range("AE1").Value = "SCT_sys"
$ws.Range("AF1").Value = "SCT_idle"
$ws.Range("AG1").Value = "SCT_iowait"
$ws.Range("AH1").Value = "SCT_irq"
$ws.Range("AI1").Value = "SCT_softirq"
$ws.Range("AJ1").Value = "SCT_steal"
$ws.Range("AK1").Value = "SCT_guest"
$ws.Range("AL1").Value = "SCT_guest_nice"
$ws.Range("AM1").Value = "Cor_in_Sys"
$ws.Range("AN1").Value = "CS_ctx_switches"
$ws.Range("AO1").Value = "CS_interrupts"
$ws.Range("AP1").Value = "CS_soft_interrupts"
$ws.Range("AQ1").Value = "CS_syscalls"
$ws.Range("AR1").Value = "SMU_total"
$ws.Range("AS1").Value = "SMU_available"
$ws.Range("AT1").Value = "SMU_percent"
$ws.Range("AU1").Value = "SMU_used"
$ws.Range("AV1").Value = "SMU_free"
$ws.Range("AW1").Value = "SMU_active"
$ws.Range("AX1").Value = "SMU_inactive"
$ws.Range("AY1").Value = "SMU_buffers"
$ws.Range("AZ1").Value = "SMU_cached"
$ws.Range("BA1").Value = "SMU_shared"
$ws.Range("BB1").Value = "SMU_slab"
$ws.Range("BC1").Value = "DU_total"
$ws.Range("BD1").Value = "DU_used"
$ws.Range("BE1").Value = "DU_free"
$ws.Range("BF1").Value = "DU_percent"
$ws.Range("BG1").Value = "NIOB_sent"
$ws.Range("BH1").Value = "NIOB_received"
$ws.Range("BI1").Value = "NIOP_sent"
$ws.Range("BJ1").Value = "NIO_received"
$ws.Range("BK1").Value = "NIO_errin"
$ws.Range("BL1").Value = "NIO_errout"
$ws.Range("BM1").Value = "NIO_dropin"
$ws.Range("BN1").Value = "NIO_dropout"

# Update the saved selection / scroll anchor to match the authored state.
$ws.Range("I16").Select()
